# Refresh the crypto price/volume table (columns D "Price" and E "Volume(1h)")
# with the latest scraped figures. Price strings that look like a plain number
# (e.g. "608.30") are written with a leading apostrophe so Excel keeps them as
# text (matching the original quote-prefixed/inline-string cells) instead of
# silently auto-converting them to numeric values and dropping the trailing
# zero formatting; multi-dot prices (e.g. "66.261.55") are never auto-numeric
# so they're assigned as plain strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.261.55'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '3.536.03'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''608.30'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').Value = '''144.15'
$ws.Range('E6').Value = '  -2.81%  '
$ws.Range('D7').Value = '3.533.80'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  -4.23%  '
$ws.Range('D11').Value = '''8.05'
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('E12').Value = '  -2.71%  '
$ws.Range('D13').Value = '4.133.78'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '''0.0000208'
$ws.Range('E14').Value = '  -4.31%  '
$ws.Range('D15').Value = '''30.39'
$ws.Range('E15').Value = '  -4.93%  '
$ws.Range('D16').Value = '3.534.38'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').Value = '66.300.69'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').Value = '''0.115'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '''10.96'
$ws.Range('E19').Value = '  +1.76%  '
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').Value = '''425.90'
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range('D23').Value = '''0.602'
$ws.Range('E23').Value = '  -1.58%  '
$ws.Range('D24').Value = '''78.79'
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').Value = '3.674.26'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('D28').Value = '''9.27'
$ws.Range('E28').Value = '  -5.50%  '
$ws.Range('D29').Value = '''8.05'
$ws.Range('E29').Value = '  -3.16%  '
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('E32').Value = '  -3.78%  '
$ws.Range('D33').Value = '''1.49'
$ws.Range('E33').Value = '  -6.49%  '
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').Value = '3.522.99'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -3.07%  '
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('D39').Value = '''5.63'
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '''172.12'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').Value = '''0.0858'
$ws.Range('E42').Value = '  -4.06%  '
$ws.Range('E43').Value = '  -5.11%  '
$ws.Range('D44').Value = '''0.894'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('E45').Value = '  -8.97%  '
$ws.Range('D46').Value = '''45.26'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').Value = '''26.04'
$ws.Range('E47').Value = '  -7.82%  '
$ws.Range('E48').Value = '  -7.61%  '
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('E50').Value = '  -4.18%  '
$ws.Range('D51').Value = '''0.951'
$ws.Range('E51').Value = '  -4.16%  '
